$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 587, pushing existing rows 587:681 down to 589:683
$ws.Rows("587:588").Insert()

# Populate new row 587
$ws.Range("A587").Value = 9
$ws.Range("B587").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C587").Value = "Metropolitana"
$ws.Range("D587").Value = "2023-01-20"
$ws.Range("E587").Value = 13
$ws.Range("F587").Value = "Fruta"
$ws.Range("G587").Value = 100109
$ws.Range("H587").Value = "Uva"
$ws.Range("I587").Value = 100109001
$ws.Range("J587").Value = "Uva"
$ws.Range("K587").Value = "Flame Seedless"
$ws.Range("L587").Value = "Primera"
$ws.Range("M587").Value = 120
$ws.Range("N587").Value = 9000
$ws.Range("O587").Value = 9000
$ws.Range("P587").Value = 9000
$ws.Range("Q587").Value = '$/bandeja 10 kilos'
$ws.Range("R587").Value = "Provincia de Limarí"
$ws.Range("S587").Value = 900
$ws.Range("T587").Value = 10

# Populate new row 588
$ws.Range("A588").Value = 9
$ws.Range("B588").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C588").Value = "Metropolitana"
$ws.Range("D588").Value = "2023-01-20"
$ws.Range("E588").Value = 13
$ws.Range("F588").Value = "Fruta"
$ws.Range("G588").Value = 100109
$ws.Range("H588").Value = "Uva"
$ws.Range("I588").Value = 100109001
$ws.Range("J588").Value = "Uva"
$ws.Range("K588").Value = "Superior Seedless"
$ws.Range("L588").Value = "Primera"
$ws.Range("M588").Value = 150
$ws.Range("N588").Value = 10000
$ws.Range("O588").Value = 10000
$ws.Range("P588").Value = 10000
$ws.Range("Q588").Value = '$/bandeja 10 kilos'
$ws.Range("R588").Value = "Provincia de Limarí"
$ws.Range("S588").Value = 1000
$ws.Range("T588").Value = 10
